# Cell Score Correlation Matrix workbook update
# Removes three gene-program rows that were dropped from the curated list:
#   - "Li (2019) Mature macrophage"            (originally row 12)
#   - "Li (2019) Terminal maturation macrophage" (originally row 16)
#   - "Han (2020) T_cell"                      (originally row 22)
# Deleting bottom-to-top keeps the earlier row numbers valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(22).Delete()
$ws.Rows(16).Delete()
$ws.Rows(12).Delete()

# Leave the selection where the author's Excel session ended up: the row
# that is now "Han (2020) Fetal epithelial progenitor" (row 15), selected
# as a full row.
$ws.Rows(15).Select()
